# Replace row 6 (PMC7091225 / "Clinical features of severe pediatric ...")
# with the content belonging to PMC7095102 ("Characteristics of pediatric
# SARS-CoV-2 infection and potential evidence for persistent fecal viral
# shedding"), matching the data already used for row 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstractText = @"
We report epidemiological and clinical investigations on ten pediatric SARS-CoV-2 infection cases confirmed by real-time reverse transcription PCR assay of SARS-CoV-2 RNA.
 Symptoms in these cases were nonspecific and no children required respiratory support or intensive care.
 Chest X-rays lacked definite signs of pneumonia, a defining feature of the infection in adult cases.
 Notably, eight children persistently tested positive on rectal swabs even after nasopharyngeal testing was negative, raising the possibility of fecal–oral transmission.

"@

$authorsText = @"
[Yi%Xu%NULL%0,    Xufang%Li%NULL%1,    Bing%Zhu%NULL%1,    Huiying%Liang%NULL%1,    Chunxiao%Fang%NULL%1,    Yu%Gong%NULL%1,    Qiaozhi%Guo%NULL%1,    Xin%Sun%NULL%1,    Danyang%Zhao%NULL%1,    Jun%Shen%NULL%1,    Huayan%Zhang%NULL%1,    Hongsheng%Liu%NULL%2,    Hongsheng%Liu%NULL%0,    Huimin%Xia%huiminxia@hotmail.com%1,    Jinling%Tang%jltang@cuhk.edu.hk%1,    Kang%Zhang%kang.zhang@gmail.com%1,    Sitang%Gong%sitangg@126.com%2,    Sitang%Gong%sitangg@126.com%0]
"@

$ws.Range("C6").Value = $ws.Range("C7").Text
$ws.Range("D6").Value = $abstractText
$ws.Range("E6").Value = $authorsText
$ws.Range("F6").Value = $ws.Range("F7").Text
